$d = $word.ActiveDocument

$rng = $d.Content.Find.Execute("blue", $true, $false, $false, $false, $false, $true, 1, $false, "orange", 2)
